$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = '$191.8 \pm 19.57$'
$ws.Range("B3").Value = '$388.4 \pm 18.84$'
$ws.Range("B4").Value = '$455.7 \pm 55.28$'
$ws.Range("B5").Value = '$536.5 \pm 35.26$'
$ws.Range("B6").Value = '$998.9 \pm 54.48$'
$ws.Range("B7").Value = '$1516 \pm 56.59$'
$ws.Range("B8").Value = '$2526 \pm 119.5$'
$ws.Range("B9").Value = '$5743 \pm 361.3$'
$ws.Range("C2").Value = '$0.4187 \pm 0.4342$'
$ws.Range("C3").Value = '$0.943 \pm 0.3481$'
$ws.Range("C4").Value = '$0.6104 \pm 0.3138$'
$ws.Range("C5").Value = '$1.743 \pm 0.9464$'
$ws.Range("C6").Value = '$5.216 \pm 0.4816$'
$ws.Range("C7").Value = '$5.609 \pm 0.8756$'
$ws.Range("C8").Value = '$11.45 \pm 3.29$'
$ws.Range("C9").Value = '$13.69 \pm 4.05$'
$ws.Range("D2").Value = '$0.034 \pm 0$'
$ws.Range("D3").Value = '$0.0904 \pm 0.0224$'
$ws.Range("D4").Value = '$0.1352 \pm 0.0364$'
$ws.Range("D5").Value = '$0.1948 \pm 0.0395$'
$ws.Range("D6").Value = '$0.8726 \pm 0.1529$'
$ws.Range("D7").Value = '$1.631 \pm 0.3915$'
$ws.Range("D8").Value = '$3.191 \pm 0.4752$'
$ws.Range("D9").Value = '$7.635 \pm 1.872$'
$ws.Range("E2").Value = '$0.0367 \pm 0.0041$'
$ws.Range("E3").Value = '$0.0705 \pm 0.0129$'
$ws.Range("E4").Value = '$0.0724 \pm 0.0051$'
$ws.Range("E5").Value = '$0.1203 \pm 0.0184$'
$ws.Range("E6").Value = '$0.9277 \pm 0.2418$'
$ws.Range("E7").Value = '$0.8938 \pm 0.1717$'
$ws.Range("E8").Value = '$1.275 \pm 0.3423$'
$ws.Range("E9").Value = '$3.817 \pm 0.5884$'
